# Rename the original sheet and add the new "Employee" sheet right after it.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Events"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Employee"

# Header row on the new sheet (bold, matches the existing header style).
$ws2.Range("A1:B1").Font.Bold = $true
$ws2.Range("A1").Value = "EmployeeID"
$ws2.Range("B1").Value = "EmployeeName"

# Update the Events sheet header from "Fjöldi" to "Employees".
$ws1.Range("C1").Value = "Employees"

# Employee rows, entered in the same order the shared-string table ends up
# in (Kata before Inga/Jóna), even though the final sheet lists them
# alphabetically by name.
$ws2.Range("B2").Value = "Anna"
$ws2.Range("B3").Value = "Brynja"
$ws2.Range("B4").Value = "Diljá"
$ws2.Range("B5").Value = "Emma"
$ws2.Range("B6").Value = "Freyja"
$ws2.Range("B7").Value = "Gunna"
$ws2.Range("B8").Value = "Hildur"
$ws2.Range("B11").Value = "Kata"
$ws2.Range("B9").Value = "Inga"
$ws2.Range("B10").Value = "Jóna"
$ws2.Range("B12").Value = "Lísa"

for ($r = 2; $r -le 12; $r++) {
    $ws2.Cells.Item($r, 1).Value = $r - 1
}

# Resize columns C and D on the Events sheet to their new best-fit widths.
# (COM's ColumnWidth setter here snaps to a 1/6-character pixel grid, so the
# inputs are pre-compensated by that fixed +0.8333... padding offset to land
# as close as possible to the true OOXML widths of 9.90625 / 8.90625.)
$ws1.Columns.Item(3).ColumnWidth = 9.072916666666666
$ws1.Columns.Item(4).ColumnWidth = 8.072916666666666

# Column widths for the new sheet (true OOXML widths 10.36328125 /
# 13.54296875, compensated the same way).
$ws2.Columns.Item(1).ColumnWidth = 9.529947916666666
$ws2.Columns.Item(2).ColumnWidth = 12.709635416666666

# Selections: Events shows C1 selected (no longer the active tab); Employee
# becomes the active tab with C10 selected.
$ws1.Range("C1").Select()
$ws2.Range("C10").Select()
$ws2.Activate()
